$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''72.399.39'
$ws.Range("E2").Value = '  -0.50%  '
$ws.Range("D3").Value = '''3.921.86'
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '''595.51'
$ws.Range("E5").Value = '  +0.78%  '
$ws.Range("D6").Value = '''169.23'
$ws.Range("E6").Value = '  +10.56%  '
$ws.Range("D7").Value = '''0.680'
$ws.Range("E7").Value = '  -0.90%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").Value = '''0.777'
$ws.Range("E9").Value = '  +2.73%  '
$ws.Range("D10").Value = '''0.186'
$ws.Range("E10").Value = '  +9.83%  '
$ws.Range("D11").Value = '''55.16'
$ws.Range("E11").Value = '  +1.45%  '
$ws.Range("D12").Value = '''0.0000328'
$ws.Range("E12").Value = '  +2.72%  '
$ws.Range("D13").Value = '''11.45'
$ws.Range("E13").Value = '  +4.42%  '
$ws.Range("D14").Value = '''4.546.15'
$ws.Range("E14").Value = '  -2.29%  '
$ws.Range("D15").Value = '''3.926.37'
$ws.Range("E15").Value = '  -2.09%  '
$ws.Range("D16").Value = '''21.23'
$ws.Range("E16").Value = '  +2.90%  '
$ws.Range("D17").Value = '''14.08'
$ws.Range("E17").Value = '  -0.76%  '
$ws.Range("D18").Value = '''1.22'
$ws.Range("E18").Value = '  -4.21%  '
$ws.Range("D19").Value = '''72.394.19'
$ws.Range("E19").Value = '  -0.44%  '
$ws.Range("E20").Value = '  -1.30%  '
$ws.Range("D21").Value = '''445.80'
$ws.Range("E21").Value = '  +1.97%  '
$ws.Range("D22").Value = '''4.76'
$ws.Range("E22").Value = '  -1.06%  '
$ws.Range("D23").Value = '''95.11'
$ws.Range("E23").Value = '  -1.72%  '
$ws.Range("D24").Value = '''3.30'
$ws.Range("E24").Value = '  -6.08%  '
$ws.Range("D25").Value = '''13.99'
$ws.Range("E25").Value = '  -2.54%  '
$ws.Range("D26").Value = '''4.24'
$ws.Range("E26").Value = '  -3.57%  '
$ws.Range("D27").Value = '''11.10'
$ws.Range("E27").Value = '  -2.76%  '
$ws.Range("D28").Value = '''5.94'
$ws.Range("E28").Value = '  +0.27%  '
$ws.Range("D29").Value = '''10.30'
$ws.Range("E29").Value = '  -4.20%  '
$ws.Range("D30").Value = '''35.59'
$ws.Range("D31").Value = '''7.86'
$ws.Range("E31").Value = '  -1.14%  '
$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").Value = '''50.71'
$ws.Range("E32").Value = '  +0.98%  '
$ws.Range("B33").Value = 'Cosmos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D33").Value = '''13.75'
$ws.Range("E33").Value = '  +0.50%  '
$ws.Range("E34").Value = '  -4.58%  '
$ws.Range("D35").Value = '''0.0₃0979'
$ws.Range("E35").Value = '  +13.29%  '
$ws.Range("D36").Value = '''68.53'
$ws.Range("E36").Value = '  -3.87%  '
$ws.Range("D37").Value = '''620.90'
$ws.Range("E37").Value = '  -9.02%  '
$ws.Range("D38").Value = '''0.424'
$ws.Range("E38").Value = '  -4.03%  '
$ws.Range("D39").Value = '''1.00'
$ws.Range("E39").Value = '  +0.14%  '
$ws.Range("D40").Value = '''3.32'
$ws.Range("E40").Value = '  -0.30%  '
$ws.Range("D41").Value = '''0.144'
$ws.Range("E41").Value = '  -2.86%  '
$ws.Range("E42").Value = '  +0.14%  '
$ws.Range("D43").Value = '''3.21'
$ws.Range("E43").Value = '  +41.14%  '
$ws.Range("D44").Value = '''0.0474'
$ws.Range("E44").Value = '  -3.57%  '
$ws.Range("D45").Value = '''10.43'
$ws.Range("E45").Value = '  -6.46%  '
$ws.Range("D46").Value = '''0.146'
$ws.Range("E46").Value = '  -2.69%  '
$ws.Range("B47").Value = 'WEMIXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").Value = '''2.88'
$ws.Range("E47").Value = '  -15.23%  '
$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D48").Value = '''3.36'
$ws.Range("E48").Value = '  -0.89%  '
$ws.Range("B49").Value = 'Fetch.AI'
$ws.Range("C49").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D49").Value = '''2.55'
$ws.Range("E49").Value = '  -8.67%  '
$ws.Range("D50").Value = '''2.829.07'
$ws.Range("E50").Value = '  -0.28%  '
$ws.Range("D51").Value = '''0.000274'
$ws.Range("E51").Value = '  +1.72%  '
